# Build site at 2022-09-26 16:07:08 UTC
# Applies the LOQ4058 content update:
#  - removes the standalone "Docentes responsaveis" value row (old row 13),
#    which shifts every subsequent row up by one
#  - fixes up the handful of cells whose text content changed / got
#    reshuffled as part of the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 13 (the standalone B13/C13 "1488970 - Marivone Nunho
# Sousa" row sitting under "Docentes responsaveis:") - this shifts rows
# 14-25 up to become rows 13-24.
$ws.Rows.Item(13).Delete()

# After the shift, patch the cells whose content ended up different from a
# pure shift (per the target workbook).
$ws.Range("B10").Value = "1488970 - Marivone Nunho Sousa"
$ws.Range("C10").Value = "1488970 - Marivone Nunho Sousa"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

$ws.Range("B18").Value = "1488970 - Marivone Nunho Sousa"
$ws.Range("C18").Value = "1488970 - Marivone Nunho Sousa"

$ws.Range("B19").Value = "Participação em sala de aula, preparação e apresentação de trabalhos e provas escritas."
$ws.Range("C19").Value = "Participação em sala de aula, preparação e apresentação de trabalhos e provas escritas."

$ws.Range("B20").Value = "Média Final = (Prova1 + Prova2 + Nota de Trabalho) /3`nMédia final mínima de aprovação = 5,0"
$ws.Range("C20").Value = "Média Final = (Prova1 + Prova2 + Nota de Trabalho) /3`nMédia final mínima de aprovação = 5,0"

$ws.Range("B21").Value = "(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0"
$ws.Range("C21").Value = "(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0"
